$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-13 Tuesday", "2026-01-14 Wednesday"),
    @("30÷3=10, 0", "44÷5=8, 4"),
    @("44÷4=11, 0", "69÷8=8, 5"),
    @("60÷7=8, 4", "52÷7=7, 3"),
    @("79÷5=15, 4", "96÷7=13, 5"),
    @("56÷5=11, 1", "21÷4=5, 1"),
    @("98÷5=19, 3", "12÷6=2, 0"),
    @("84÷3=28, 0", "49÷8=6, 1"),
    @("97÷2=48, 1", "45÷7=6, 3"),
    @("70÷2=35, 0", "96÷4=24, 0"),
    @("80÷8=10, 0", "41÷4=10, 1"),
    @("81÷2=40, 1", "35÷8=4, 3"),
    @("48÷9=5, 3", "86÷8=10, 6"),
    @("33÷8=4, 1", "19÷5=3, 4"),
    @("66÷2=33, 0", "54÷5=10, 4"),
    @("35÷3=11, 2", "78÷5=15, 3"),
    @("57÷2=28, 1", "87÷7=12, 3"),
    @("20÷7=2, 6", "98÷8=12, 2"),
    @("51÷6=8, 3", "81÷5=16, 1"),
    @("40÷9=4, 4", "23÷8=2, 7"),
    @("93÷5=18, 3", "88÷9=9, 7"),
    @("80÷2=40, 0", "51÷2=25, 1"),
    @("34÷5=6, 4", "34÷9=3, 7"),
    @("24÷9=2, 6", "83÷4=20, 3"),
    @("36÷5=7, 1", "63÷8=7, 7"),
    @("96÷6=16, 0", "33÷5=6, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
